$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "Ankita Gangotra"

$ws.Range("A25").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = "A.G"

$ws.Range("D25").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = [DateTime]"2014-06-10"

$ws.Range("D31").Select()
